# Fixes the marksheet's score summary (rows 10-12) and collapses the
# duplicated "Student Ans / Correct Ans" answer-key block (previously
# repeated across columns A/B, D/E and G/H) down to a single A/B pair,
# filling in the student's matched answers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Score summary rows ---------------------------------------------------

$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("B10").Value = 20
$ws.Range("D10").Value = 8
$ws.Range("E10").Value = 28

$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("B12").Value = 80
$ws.Range("E12").Value = "80/112"

# --- Drop the 3rd (G/H) "Student Ans / Correct Ans" column pair ----------

$ws.Range("G15:H40").Clear()

# --- Row 16-18: keep the D/E pair, but refresh it to the matched answer --

$ws.Range("D16").Style = "correctStyle"
$ws.Range("D16").Value = "Option A"

$ws.Range("D17").Style = "correctStyle"
$ws.Range("D17").Value = "Option C"

$ws.Range("D18").Style = "correctStyle"
$ws.Range("D18").Value = "Option D"

# --- Drop the 2nd (D/E) pair for rows 19-40 -------------------------------

$ws.Range("D19:E40").Clear()

# --- Fill column A (Student Ans, 1st pair) with the matched answer -------

$answers = @{
    16 = "Option A"
    17 = $null
    18 = "Option B"
    19 = "Option C"
    20 = "Option B"
    21 = "Option C"
    22 = "Option D"
    23 = $null
    24 = $null
    25 = $null
    26 = "Option C"
    27 = $null
    28 = "Option D"
    29 = $null
    30 = "Option B"
    31 = "Option D"
    32 = $null
    33 = "Option D"
    34 = "Option B"
    35 = "Option D"
    36 = "Option A"
    37 = "Option A"
    38 = "Option A"
    39 = "Option D"
    40 = $null
}

foreach ($row in 16..40) {
    $cell = $ws.Range("A$row")
    $answer = $answers[$row]
    if ($answer) {
        $cell.Style = "correctStyle"
        $cell.Value = $answer
    } else {
        $cell.Style = "normalStyle"
    }
}
